$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.38
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 4.8
$ws.Range("K2").Value = 5.2
$ws.Range("L2").Value = 1.51
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 1.74
$ws.Range("Q2").Value = 2.28
$ws.Range("R2").Value = 1.26
$ws.Range("S2").Value = 4.5
$ws.Range("T2").Value = 2.76
$ws.Range("V2").Value = 1.08
$ws.Range("W2").Value = 3.5
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 130
$ws.Range("AB2").Value = 5.7
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 60
$ws.Range("AE2").Value = 610
$ws.Range("AF2").Value = 6.4
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 55
$ws.Range("AI2").Value = 420
$ws.Range("AJ2").Value = 11
$ws.Range("AK2").Value = 19.5
$ws.Range("AL2").Value = 75
$ws.Range("AN2").Value = 10
$ws.Range("AO2").Value = 850
$ws.Range("F3").Value = 1.32
$ws.Range("G3").Value = 1.38
$ws.Range("I3").Value = 10.5
$ws.Range("J3").Value = 5.9
$ws.Range("K3").Value = 7.4
$ws.Range("L3").Value = 1.24
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 7.4
$ws.Range("O3").Value = 1.13
$ws.Range("P3").Value = 3.2
$ws.Range("Q3").Value = 1.41
$ws.Range("R3").Value = 1.89
$ws.Range("S3").Value = 2.08
$ws.Range("T3").Value = 1.69
$ws.Range("U3").Value = 2.12
$ws.Range("W3").Value = 3.55
$ws.Range("AN3").Value = 5.8
$ws.Range("F4").Value = 1.9
$ws.Range("G4").Value = 2.18
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 4.7
$ws.Range("J4").Value = 3.4
$ws.Range("K4").Value = 4.2
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 4
$ws.Range("O4").Value = 1.27
$ws.Range("P4").Value = 2.08
$ws.Range("Q4").Value = 1.83
$ws.Range("R4").Value = 1.44
$ws.Range("S4").Value = 3.05
$ws.Range("T4").Value = 1.7
$ws.Range("U4").Value = 2.06
$ws.Range("V4").Value = 1.28
$ws.Range("W4").Value = 1.85
$ws.Range("X4").Value = 17.5
$ws.Range("Y4").Value = 17.5
$ws.Range("Z4").Value = 34
$ws.Range("AA4").Value = 290
$ws.Range("AB4").Value = 10.5
$ws.Range("AC4").Value = 9.199999999999999
$ws.Range("AD4").Value = 18
$ws.Range("AE4").Value = 55
$ws.Range("AF4").Value = 14
$ws.Range("AG4").Value = 11.5
$ws.Range("AH4").Value = 19
$ws.Range("AI4").Value = 60
$ws.Range("AJ4").Value = 26
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 38
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 15
$ws.Range("AO4").Value = 55
$ws.Range("H5").Value = 1.5
$ws.Range("I5").Value = 1.51
$ws.Range("J5").Value = 4.7
$ws.Range("K5").Value = 4.8
$ws.Range("S5").Value = 3.35
$ws.Range("T5").Value = 2.06
$ws.Range("U5").Value = 1.9
$ws.Range("V5").Value = 2.96
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 7.6
$ws.Range("Z5").Value = 8
$ws.Range("AE5").Value = 16
$ws.Range("AG5").Value = 29
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 38
$ws.Range("AJ5").Value = 270
$ws.Range("AK5").Value = 130
$ws.Range("AM5").Value = 150
$ws.Range("AN5").Value = 180
$ws.Range("AO5").Value = 8
$ws.Range("F6").Value = 1.9
$ws.Range("G6").Value = 2.24
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 4.5
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 4.3
$ws.Range("N6").Value = 3.6
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 1.96
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 1.38
$ws.Range("S6").Value = 2.8
$ws.Range("T6").Value = 1.67
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.29
$ws.Range("W6").Value = 1.81
$ws.Range("AN6").Value = 65
$ws.Range("F7").Value = 1.8
$ws.Range("G7").Value = 1.92
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 5.8
$ws.Range("J7").Value = 3.55
$ws.Range("L7").Value = 1.44
$ws.Range("N7").Value = 3.4
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 1.84
$ws.Range("Q7").Value = 1.98
$ws.Range("R7").Value = 1.31
$ws.Range("S7").Value = 3.5
$ws.Range("T7").Value = 1.88
$ws.Range("U7").Value = 1.89
$ws.Range("V7").Value = 1.21
$ws.Range("W7").Value = 2.08
$ws.Range("X7").Value = 14.5
$ws.Range("Y7").Value = 17.5
$ws.Range("Z7").Value = 38
$ws.Range("AA7").Value = 130
$ws.Range("AB7").Value = 8.4
$ws.Range("AC7").Value = 9.4
$ws.Range("AD7").Value = 21
$ws.Range("AE7").Value = 75
$ws.Range("AF7").Value = 11
$ws.Range("AH7").Value = 22
$ws.Range("AI7").Value = 85
$ws.Range("AJ7").Value = 21
$ws.Range("AK7").Value = 21
$ws.Range("AN7").Value = 14.5
$ws.Range("AO7").Value = 90
